$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label (A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 11:20"

# Row 25: Chequia
$ws.Range("A25").Value = "Chequia"
$ws.Range("B25").Value = 2859
$ws.Range("C25").Value = 42
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 2831
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 17

# Row 26: Malasia
$ws.Range("A26").Value = "Malasia"
$ws.Range("B26").Value = 2626
$ws.Range("C26").Value = 156
$ws.Range("D26").Value = 388
$ws.Range("E26").Value = 2201
$ws.Range("F26").Value = 94
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 37

# Row 27: Irlanda
$ws.Range("A27").Value = "Irlanda"
$ws.Range("B27").Value = 2615
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 2564
$ws.Range("F27").Value = 59
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 46

# Row 28: Dinamarca
$ws.Range("A28").Value = "Dinamarca"
$ws.Range("B28").Value = 2555
$ws.Range("C28").Value = 160
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 2482
$ws.Range("F28").Value = 113
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 72

# Row 29: Chile
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 2139
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 75
$ws.Range("E29").Value = 2057
$ws.Range("F29").Value = 7
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 7

# Row 39: Indonesia
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 1414
$ws.Range("C39").Value = 129
$ws.Range("D39").Value = 75
$ws.Range("E39").Value = 1217
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 8
$ws.Range("H39").Value = 122

# Row 40: Arabia Saudita
$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("B40").Value = 1299
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 66
$ws.Range("E40").Value = 1225
$ws.Range("F40").Value = 12
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 8

# Row 41: Finlandia
$ws.Range("A41").Value = "Finlandia"
$ws.Range("B41").Value = 1286
$ws.Range("C41").Value = 46
$ws.Range("D41").Value = 10
$ws.Range("E41").Value = 1265
$ws.Range("F41").Value = 32
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 11

# Row 81: Kazajistan
$ws.Range("A81").Value = "Kazajistan"
$ws.Range("B81").Value = 294
$ws.Range("C81").Value = 10
$ws.Range("D81").Value = 20
$ws.Range("E81").Value = 273
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1

# Row 87: Albania
$ws.Range("A87").Value = "Albania"
$ws.Range("B87").Value = 223
$ws.Range("C87").Value = 11
$ws.Range("D87").Value = 44
$ws.Range("E87").Value = 168
$ws.Range("F87").Value = 7
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 11

# Row 88: Burkina Faso
$ws.Range("A88").Value = "Burkina Faso"
$ws.Range("B88").Value = 222
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 23
$ws.Range("E88").Value = 187
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 12

# Row 89: Republica de Chipre
$ws.Range("A89").Value = "Republica de Chipre"
$ws.Range("B89").Value = 214
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 15
$ws.Range("E89").Value = 193
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 6

# Row 127: Isla de Man
$ws.Range("A127").Value = "Isla de Man"
$ws.Range("B127").Value = 46
$ws.Range("C127").Value = 4
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 46
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0

# Row 128: Monaco
$ws.Range("A128").Value = "Monaco"
$ws.Range("B128").Value = 46
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 44
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 1

# Row 129: Guayana Francesa
$ws.Range("A129").Value = "Guayana Francesa"
$ws.Range("B129").Value = 43
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 6
$ws.Range("E129").Value = 37
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0
